$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab) to match the new name.
$ws.Name = "defaultDialog.csv"

# Append the new dialog rows (36-41), each following the same pattern as the
# preceding rows: column A is the dialog id string, B is 4, C/D/E/F are 0.
$newRows = @(
    @{ Row = 36; Text = "dialog_equip_an_equipment_success" },
    @{ Row = 37; Text = "dialog_unequip_a_shipheader" },
    @{ Row = 38; Text = "dialog_equip_a_shipheader" },
    @{ Row = 39; Text = "dialog_cannot_unequip_a_shipheader_demon" },
    @{ Row = 40; Text = "dialog_cannot_unequip_a_shipheader_demon_first" },
    @{ Row = 41; Text = "dialog_cannot_unequip_a_shipheader" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Text
    $ws.Range("B$row").Value = 4
    $ws.Range("C$row").Value = 0
    $ws.Range("D$row").Value = 0
    $ws.Range("E$row").Value = 0
    $ws.Range("F$row").Value = 0
}

# Update the selection to match the new last row, mirroring the saved view state.
$ws.Range("B41:F41").Select()
